# Update Work Week and Social Spending
# Updates the "Data" sheet GDP per Capita values (column E) for existing rows
# (years 1820-2010) and appends 6 new rows for years 2011-2016.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update existing rows ($E$2:$E$192) ---
$dataUpdates = @{
    2 = "1117"
    32 = "867"
    42 = "856"
    52 = "845"
    62 = "864"
    72 = "881"
    82 = "926"
    95 = "974"
    107 = "1078"
    111 = "1116"
    115 = "1154"
    120 = "1549"
    124 = "1546"
    125 = "1833"
    128 = "2628"
    129 = "2410"
    132 = "2115"
    133 = "2251"
    134 = "2397"
    135 = "2695"
    136 = "2962"
    137 = "3220"
    138 = "3491"
    139 = "3934"
    140 = "3918"
    141 = "4050"
    142 = "4230"
    143 = "4307"
    144 = "4339"
    145 = "4395"
    146 = "4629"
    147 = "4894"
    148 = "4988"
    149 = "5066"
    150 = "5235"
    151 = "5547"
    152 = "6135"
    153 = "6062"
    154 = "6150"
    155 = "6583"
    156 = "6229"
    157 = "6129"
    158 = "5681"
    159 = "5501"
    160 = "5482"
    161 = "5318"
    162 = "4975"
    163 = "5040"
    164 = "5019"
    165 = "5117"
    166 = "5059"
    167 = "4814"
    168 = "4884"
    169 = "5254"
    170 = "5397"
    171 = "5748"
    172 = "6035"
    173 = "6116.58331807709"
    174 = "6292.02885816339"
    175 = "6414.47499695309"
    176 = "6517.0515896976"
    177 = "6656.28322158056"
    178 = "6648.93315693184"
    179 = "6523.63636193961"
    180 = "6432.6529564669"
    181 = "6485.32282006973"
    182 = "6531.28475813797"
    183 = "6615.19295637707"
    184 = "6660.76302247247"
    185 = "6908.66675257455"
    186 = "7004.98453646445"
    187 = "7073.77027017016"
    188 = "7286.96668258715"
    189 = "7402.17840857208"
    190 = "7353.10712687864"
    191 = "7115.37307299018"
    192 = "7025.93017720906"
}

foreach ($row in $dataUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $dataUpdates[$row]
}

# --- Append new rows for years 2011-2016 (rows 193-198) ---
$newRows = @(
    @{ Row = 193; Year = 2011.0; Value = "7141" }
    @{ Row = 194; Year = 2012.0; Value = "7053" }
    @{ Row = 195; Year = 2013.0; Value = "7017" }
    @{ Row = 196; Year = 2014.0; Value = "7005" }
    @{ Row = 197; Year = 2015.0; Value = "7024" }
    @{ Row = 198; Year = 2016.0; Value = "7084" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = 388.0
    $ws.Cells.Item($r, 2).Value = "Jamaica"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $entry.Year
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $entry.Value
}

